# InvestmentCalc.xlsx: extend the model from a 5-year to a 10-year
# project lifetime (columns H:L added after existing columns B:G),
# rescale the yearly cash-flow assumptions, move the "Residual" /
# "restricted Equity" one-off entries from year 5 to year 10, and
# recompute Present Value / Accumulated Present Value / Net Present
# Value (discount rate in B14 = 8.4%) for the new 10-year series.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newCols = @("H", "I", "J", "K", "L")

# --- Row 1: header label + extend the year index (0..10) out to column L ---
$ws.Range("A1").Value = "Project lifetime"
foreach ($c in $newCols) { $ws.Range("G1").Copy($ws.Range($c + "1")) }
$ws.Range("H1").Value = 6
$ws.Range("I1").Value = 7
$ws.Range("J1").Value = 8
$ws.Range("K1").Value = 9
$ws.Range("L1").Value = 10

# --- Row 3: Depreciation, now a flat 60000 across C1:L1 (was 120000 across C:G) ---
foreach ($c in $newCols) { $ws.Range("G3").Copy($ws.Range($c + "3")) }
$ws.Range("C3").Value = 60000
$ws.Range("D3").Value = 60000
$ws.Range("E3").Value = 60000
$ws.Range("F3").Value = 60000
$ws.Range("G3").Value = 60000
$ws.Range("H3").Value = 60000
$ws.Range("I3").Value = 60000
$ws.Range("J3").Value = 60000
$ws.Range("K3").Value = 60000
$ws.Range("L3").Value = 60000

# --- Row 4: Incoming Payments, now a flat 350000 (was 1050000 across C:G) ---
foreach ($c in $newCols) { $ws.Range("G4").Copy($ws.Range($c + "4")) }
$ws.Range("C4").Value = 350000
$ws.Range("D4").Value = 350000
$ws.Range("E4").Value = 350000
$ws.Range("F4").Value = 350000
$ws.Range("G4").Value = 350000
$ws.Range("H4").Value = 350000
$ws.Range("I4").Value = 350000
$ws.Range("J4").Value = 350000
$ws.Range("K4").Value = 350000
$ws.Range("L4").Value = 350000

# --- Row 5: Outgoing Payments; B5 doubles to -280000, C5:L5 flattens to -140000 ---
foreach ($c in $newCols) { $ws.Range("G5").Copy($ws.Range($c + "5")) }
$ws.Range("B5").Value = -280000
$ws.Range("C5").Value = -140000
$ws.Range("D5").Value = -140000
$ws.Range("E5").Value = -140000
$ws.Range("F5").Value = -140000
$ws.Range("G5").Value = -140000
$ws.Range("H5").Value = -140000
$ws.Range("I5").Value = -140000
$ws.Range("J5").Value = -140000
$ws.Range("K5").Value = -140000
$ws.Range("L5").Value = -140000

# --- Row 6: Residual one-off payment moves from year 5 (G6) to year 10 (L6) ---
$ws.Range("G6").Copy($ws.Range("L6"))
$ws.Range("G6").Clear()

# --- Row 7: restricted Equity payback moves from year 5 (G7) to year 10 (L7) ---
$ws.Range("G7").Copy($ws.Range("L7"))
$ws.Range("G7").Clear()

# --- Row 8: Yearly Net ---
foreach ($c in $newCols) { $ws.Range("G8").Copy($ws.Range($c + "8")) }
$ws.Range("B8").Value = -2480000
$ws.Range("C8").Value = 270000
$ws.Range("D8").Value = 270000
$ws.Range("E8").Value = 270000
$ws.Range("F8").Value = 270000
$ws.Range("G8").Value = 270000
$ws.Range("H8").Value = 270000
$ws.Range("I8").Value = 270000
$ws.Range("J8").Value = 270000
$ws.Range("K8").Value = 270000
$ws.Range("L8").Value = 540000

# --- Row 9: Present Value = Yearly Net discounted at the after-tax rate (B14) ---
foreach ($c in $newCols) { $ws.Range("G9").Copy($ws.Range($c + "9")) }
$ws.Range("B9").Value = -2480000
$ws.Range("C9").Value = 249077.49077490772
$ws.Range("D9").Value = 229776.2830026824
$ws.Range("E9").Value = 211970.74077738225
$ws.Range("F9").Value = 195544.9638167733
$ws.Range("G9").Value = 180392.03304130374
$ws.Range("H9").Value = 166413.31461374884
$ws.Range("I9").Value = 153517.81790936238
$ws.Range("J9").Value = 141621.6032374192
$ws.Range("K9").Value = 130647.23545887378
$ws.Range("L9").Value = 241046.55988722097

# --- Row 10: Accumulated Present Value (running total of row 9) ---
# F10/G10 also pick up B10's style (6), replacing their old style (7).
foreach ($c in $newCols) { $ws.Range("B10").Copy($ws.Range($c + "10")) }
$ws.Range("B10").Copy($ws.Range("F10"))
$ws.Range("B10").Copy($ws.Range("G10"))
$ws.Range("B10").Value = -2480000
$ws.Range("C10").Value = -2230922.5092250924
$ws.Range("D10").Value = -2001146.2262224099
$ws.Range("E10").Value = -1789175.4854450277
$ws.Range("F10").Value = -1593630.5216282543
$ws.Range("G10").Value = -1413238.4885869506
$ws.Range("H10").Value = -1246825.1739732018
$ws.Range("I10").Value = -1093307.3560638395
$ws.Range("J10").Value = -951685.7528264204
$ws.Range("K10").Value = -821038.5173675466
$ws.Range("L10").Value = -579991.9574803256

# --- Row 11: Net Present Value = final Accumulated Present Value; style 9 -> 8 ---
$ws.Range("B8").Copy($ws.Range("B11"))
$ws.Range("B11").Value = -579991.9574803256
